$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update column A (OrganizationalPattern names) for existing rows 2-14.
#    The old pattern names are replaced with a new set of names; rows 10-14
#    shift down one slot in the (new) shared-string table because five new
#    names were inserted ahead of them, which is handled automatically by
#    just writing the new text values here.
# ---------------------------------------------------------------------------
$newNames = @(
    "Apprenticeship",
    "ArchitectAlsoImplements",
    "ArchitectControlsProduct",
    "ArchitectureTeam",
    "CodeOwnership",
    "DeployAlongTheGrain",
    "DeveloperControlsProcess",
    "DevelopingInPairs",
    "DevelopingInPairs",
    "DistributeWorkEvenly",
    "DivideAndConquer",
    "DomainExpertiseInRoles",
    "FeatureAssignment"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newNames[$i]
}

# ---------------------------------------------------------------------------
# 2) Row 2: swap the G2/J2 values (0.89 <-> 1.01).
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 1.01
$ws.Range("J2").Value = 0.89

# ---------------------------------------------------------------------------
# 3) Append five brand-new rows (15-19), each following the same B:L
#    (1.01/0.89 trigram-frequency) pattern and an incrementing PatternIndex
#    in column M. Insert by copying row 14 downward so the cell styling
#    (borders/alignment/font) carries over exactly, then overwrite values.
# ---------------------------------------------------------------------------
$newRows = @(
    "Stand-UpMeeting",
    "StandardsLinkingLocations",
    "LockEmUpTogether",
    "LooseInterfaces",
    "OrganizationFollowsMarket"
)

$lastDataRow = 14
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $targetRow = $lastDataRow + 1 + $i
    $srcRow = $targetRow - 1

    $ws.Rows($srcRow).Copy()
    $ws.Rows($targetRow).Insert()

    # The row-insert sometimes drops the thin border that column A needs;
    # reapply it explicitly so it matches the other pattern-name cells.
    $ws.Range("A" + $targetRow).Borders.LineStyle = 1

    $ws.Cells.Item($targetRow, 1).Value = $newRows[$i]
    $ws.Range("B" + $targetRow + ":F" + $targetRow).Value = 1.01
    $ws.Range("G" + $targetRow + ":I" + $targetRow).Value = 0.89
    $ws.Range("J" + $targetRow + ":L" + $targetRow).Value = 1.01
    $ws.Cells.Item($targetRow, 13).Value = $lastDataRow + $i
}

# ---------------------------------------------------------------------------
# 4) Selection: the saved view now has M19 selected (single cell).
# ---------------------------------------------------------------------------
$ws.Range("M19").Select()
